$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Update the report title (re-run date 2025-11-19 -> 2025-11-20)
# ------------------------------------------------------------------
$ws.Range('A1').Value = '萊爾富 工作統計表  篩選月份：202511   (  製表日期:2025-11-20  )'

# ------------------------------------------------------------------
# 2) Insert two blank rows before the old row 110, pushing the
#    existing rows 110-111 down to 112-113. Copy banding/format from
#    the two rows directly above (108-109) so the new rows keep the
#    same alternating fill/border styling as the rest of the table.
# ------------------------------------------------------------------
$ws.Range('A108:AK109').Copy()
$ws.Range('A110:AK111').Insert(-4121)

$ws.Range('A108:AK109').Copy()
$ws.Range('A110:AK111').PasteSpecial(-4122)

# ------------------------------------------------------------------
# 3) Append four more blank rows (114-117) at the bottom of the
#    table, again copying the banding format from rows 108-109.
# ------------------------------------------------------------------
$ws.Range('A108:AK109').Copy()
$ws.Range('A114:AK115').PasteSpecial(-4122)

$ws.Range('A108:AK109').Copy()
$ws.Range('A116:AK117').PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 4) Fill in the new / shifted data rows.
# ------------------------------------------------------------------
$ws.Range('A110').Value = 108
$ws.Range('B110').Value = '維修'
$ws.Range('C110').Value = 2025112573
$ws.Range('D110').Value = 'E4917114111901'
$ws.Range('E110').Value = '一般件'
$ws.Range('F110').Value = 4917
$ws.Range('G110').Value = '板橋翠華店'
$ws.Range('H110').Value = '新北市板橋區'
$ws.Range('I110').Value = '2025-11-19 11:04:25'
$ws.Range('J110').Value = '星期三'
$ws.Range('K110').Value = '上午'
$ws.Range('L110').Value = 'HLM3'
$ws.Range('M110').Value = 'HL-LIFE-ET 標籤印表機'
$ws.Range('N110').Value = 'M303'
$ws.Range('O110').Value = '無反應，不會轉動'
$ws.Range('P110').Value = '門市反應LIFE ET標籤印表機無反應，設備亮紅燈，已嘗試重啟電源並重新安裝紙捲仍異常..請台芝到店協助(無反應)'
$ws.Range('Q110').Value = 'THILF04917'
$ws.Range('R110').Value = '新北一'
$ws.Range('S110').Value = '狄澤洋'
$ws.Range('T110').Value = 1
$ws.Range('U110').Value = '已完工'
$ws.Range('V110').Value = '2025-11-19 11:16:08'
$ws.Range('W110').Value = '2025-11-20 12:00:00'
$ws.Range('X110').Value = '2025-11-20 12:45:00'
$ws.Range('Y110').Value = '2025-11-20 15:16:00'
$ws.Range('Z110').Value = 0.8
$ws.Range('AB110').Value = '到場處理'
$ws.Range('AC110').Value = '更換標籤機
換上：8187001031
換下：8187000915'
$ws.Range('A112').Value = 110
$ws.Range('B112').Value = '維修'
$ws.Range('C112').Value = 2025112607
$ws.Range('D112').Value = '14098114111901'
$ws.Range('E112').Value = '一般件'
$ws.Range('F112').Value = 4098
$ws.Range('G112').Value = '三重仁美店'
$ws.Range('H112').Value = '新北市三重區'
$ws.Range('I112').Value = '2025-11-19 14:24:06'
$ws.Range('J112').Value = '星期三'
$ws.Range('K112').Value = '下午'
$ws.Range('L112').Value = 'HL24'
$ws.Range('M112').Value = 'HL-SC主機'
$ws.Range('N112').Value = 2405
$ws.Range('O112').Value = '檔案損毀(更換硬碟)'
$ws.Range('P112').Value = 'SC(SHUTTLE6S)2025/11/19 (週三) 下午 02:20 總公司明翰來信:因4098 三重仁美店 SC第二顆硬碟發生錯誤訊息，請協助一般派工，更換門市SC第二顆硬碟，資料不備份，謝謝。...請台芝到店協助PS.若因更換HD.請跟店長宣達:1.請門市先回報代收會計 2.請確認SC的代收資料是否正確 (須與代收單據逐一核對)																											與門市確認帳務做到11/18，與通訊嘉芳確認有收到11/18的銷售				'
$ws.Range('Q112').Value = 'THILF04098'
$ws.Range('R112').Value = '新北一'
$ws.Range('S112').Value = '吳宗鴻'
$ws.Range('T112').Value = 1
$ws.Range('U112').Value = '已完工'
$ws.Range('V112').Value = '2025-11-19 14:30:22'
$ws.Range('W112').Value = '2025-11-20 10:30:00'
$ws.Range('X112').Value = '2025-11-20 12:30:00'
$ws.Range('Y112').Value = '2025-11-20 18:30:00'
$ws.Range('Z112').Value = 2
$ws.Range('AB112').Value = '到場處理'
$ws.Range('AC112').Value = '更換第二顆硬碟不備份還原完成'
$ws.Range('AK112').Value = 'O'
$ws.Range('A114').Value = 112
$ws.Range('B114').Value = '維修'
$ws.Range('C114').Value = 2025112647
$ws.Range('D114').Value = '14388114111902'
$ws.Range('E114').Value = '一般件'
$ws.Range('F114').Value = 4388
$ws.Range('G114').Value = '三重薔薇店'
$ws.Range('H114').Value = '新北市三重區'
$ws.Range('I114').Value = '2025-11-19 20:37:02'
$ws.Range('J114').Value = '星期三'
$ws.Range('K114').Value = '夜間'
$ws.Range('L114').Value = 'HLF2'
$ws.Range('M114').Value = 'HL-CCD掃描器(TM)'
$ws.Range('N114').Value = 'F201'
$ws.Range('O114').Value = '掃描無反應或感應不良'
$ws.Range('P114').Value = '門市反應TM1 CCD掃描器(HC56II TR、HC76 TR)刷讀所有商品與條碼都感應不良，有亮燈沒有逼聲+游標有在輸入的位置，已多次嘗試執行掃槍校正但都過2天又開始不好刷讀，最近是11/17與11/19下午執行掃槍校正後仍異常...需請台芝到店協助'
$ws.Range('Q114').Value = 'THILF04388'
$ws.Range('R114').Value = '新北一'
$ws.Range('S114').Value = '吳宗鴻'
$ws.Range('T114').Value = 1
$ws.Range('U114').Value = '已完工'
$ws.Range('V114').Value = '2025-11-19 20:41:31'
$ws.Range('W114').Value = '2025-11-20 14:23:00'
$ws.Range('X114').Value = '2025-11-20 14:53:00'
$ws.Range('Y114').Value = '2025-11-21 00:41:00'
$ws.Range('Z114').Value = 0.5
$ws.Range('AB114').Value = '到場處理'
$ws.Range('AC114').Value = '更換掃描槍
換下8119008847
換上8119013254'
$ws.Range('AK114').Value = 'O'
$ws.Range('A115').Value = 113
$ws.Range('B115').Value = '服務'
$ws.Range('C115').Value = 2025112704
$ws.Range('F115').Value = 4917
$ws.Range('G115').Value = '板橋翠華店'
$ws.Range('H115').Value = '新北市板橋區'
$ws.Range('Q115').Value = 'THILF04917'
$ws.Range('R115').Value = '新北一'
$ws.Range('S115').Value = '狄澤洋'
$ws.Range('T115').Value = 1
$ws.Range('U115').Value = '已完工'
$ws.Range('V115').Value = '2025-11-20 12:47:08'
$ws.Range('W115').Value = '2025-11-20 12:00:00'
$ws.Range('X115').Value = '2025-11-20 12:45:00'
$ws.Range('Z115').Value = 0.8
$ws.Range('AB115').Value = '到場處理'
$ws.Range('AC115').Value = 'PMQ4+L90'
$ws.Range('AD115').Value = 'O'
$ws.Range('AK115').Value = 'O'
$ws.Range('A116').Value = 114
$ws.Range('B116').Value = '服務'
$ws.Range('C116').Value = 2025112707
$ws.Range('F116').Value = 4098
$ws.Range('G116').Value = '三重仁美店'
$ws.Range('H116').Value = '新北市三重區'
$ws.Range('Q116').Value = 'THILF04098'
$ws.Range('R116').Value = '新北一'
$ws.Range('S116').Value = '吳宗鴻'
$ws.Range('T116').Value = 1
$ws.Range('U116').Value = '已完工'
$ws.Range('V116').Value = '2025-11-20 13:17:47'
$ws.Range('W116').Value = '2025-11-20 12:40:00'
$ws.Range('X116').Value = '2025-11-20 13:10:00'
$ws.Range('Z116').Value = 0.5
$ws.Range('AB116').Value = '到場處理'
$ws.Range('AC116').Value = 'PMQ4+L90'
$ws.Range('AD116').Value = 'O'
$ws.Range('AK116').Value = 'O'
$ws.Range('A117').Value = 115
$ws.Range('B117').Value = '服務'
$ws.Range('C117').Value = 2025112727
$ws.Range('F117').Value = 4388
$ws.Range('G117').Value = '三重薔薇店'
$ws.Range('H117').Value = '新北市三重區'
$ws.Range('Q117').Value = 'THILF04388'
$ws.Range('R117').Value = '新北一'
$ws.Range('S117').Value = '吳宗鴻'
$ws.Range('T117').Value = 1
$ws.Range('U117').Value = '已完工'
$ws.Range('V117').Value = '2025-11-20 14:55:58'
$ws.Range('W117').Value = '2025-11-20 14:05:00'
$ws.Range('X117').Value = '2025-11-20 14:15:00'
$ws.Range('Z117').Value = 0.2
$ws.Range('AB117').Value = '到場處理'
$ws.Range('AC117').Value = '+L90'
$ws.Range('AJ117').Value = 'O'
$ws.Range('AK117').Value = 'O'

# ------------------------------------------------------------------
# 5) Update the print area to cover the new extent of the table and
#    move the selection to the new last cell, matching the sheet
#    view state after the edit.
# ------------------------------------------------------------------
$ws.PageSetup.PrintArea = '$A$1:$AK$117'
$ws.Range('A117').Select()
